# Updates cryptos list values (Price column D, Volume(1h) column E)
# per the commit diff. Values are plain text (inline/shared strings) in
# the workbook, so numeric-looking Price values must be written with a
# Text number format to keep Excel from auto-converting them to numbers
# (which would also drop formatting like trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "246.26", "0.0719")
# -- force Text format first so Excel stores them as strings, not numbers.
$numericLookingCells = @("D5","D8","D10","D11","D14","D15","D19","D21","D22","D23","D24","D25","D26","D27","D29","D33","D37","D42","D43","D44","D46","D48","D50")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# New cell values (row order matches the sheet; D = Price, E = Volume(1h))
$ws.Range('D2').Value = '35.396.74'
$ws.Range('E2').Value = '  -0.20%  '

$ws.Range('D3').Value = '1.914.14'
$ws.Range('E3').Value = '  +1.79%  '

$ws.Range('D5').Value = '246.26'
$ws.Range('E5').Value = '  +2.42%  '

$ws.Range('E6').Value = '  +6.04%  '

$ws.Range('E7').Value = '  -0.36%  '

$ws.Range('D8').Value = '41.78'
$ws.Range('E8').Value = '  -2.78%  '

$ws.Range('E9').Value = '  +4.06%  '

$ws.Range('D10').Value = '53.21'
$ws.Range('E10').Value = '  +13.59%  '

$ws.Range('D11').Value = '0.0719'
$ws.Range('E11').Value = '  +2.74%  '

$ws.Range('E12').Value = '  +0.05%  '

$ws.Range('D13').Value = '2.191.75'
$ws.Range('E13').Value = '  +1.88%  '

$ws.Range('D14').Value = '12.30'
$ws.Range('E14').Value = '  +5.28%  '

$ws.Range('D15').Value = '0.701'
$ws.Range('E15').Value = '  +1.87%  '

$ws.Range('D16').Value = '1.921.77'
$ws.Range('E16').Value = '  +1.48%  '

$ws.Range('E17').Value = '  +1.36%  '

$ws.Range('D18').Value = '35.452.19'
$ws.Range('E18').Value = '  -0.03%  '

$ws.Range('D19').Value = '72.14'
$ws.Range('E19').Value = '  +2.04%  '

$ws.Range('D20').Value = '0.0₃0819'
$ws.Range('E20').Value = '  +1.98%  '

$ws.Range('D21').Value = '241.43'
$ws.Range('E21').Value = '  -0.56%  '

$ws.Range('D22').Value = '12.50'
$ws.Range('E22').Value = '  +0.67%  '

$ws.Range('D23').Value = '4.83'
$ws.Range('E23').Value = '  +1.19%  '

$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.49%  '

$ws.Range('D25').Value = '2.42'
$ws.Range('E25').Value = '  +26.36%  '

$ws.Range('D26').Value = '2.27'
$ws.Range('E26').Value = '  +0.22%  '

$ws.Range('D27').Value = '170.54'
$ws.Range('E27').Value = '  +0.09%  '

$ws.Range('E28').Value = '  +1.70%  '

$ws.Range('D29').Value = '18.41'
$ws.Range('E29').Value = '  +2.86%  '

$ws.Range('E30').Value = '  +1.19%  '

$ws.Range('D31').Value = '4.163.22'
$ws.Range('E31').Value = '  +21.94%  '

$ws.Range('E32').Value = '  +2.43%  '

$ws.Range('D33').Value = '0.0566'
$ws.Range('E33').Value = '  +0.24%  '

$ws.Range('E34').Value = '  +14.66%  '

$ws.Range('E35').Value = '  -0.54%  '

$ws.Range('E36').Value = '  +0.74%  '

$ws.Range('D37').Value = '1.75'
$ws.Range('E37').Value = '  -3.61%  '

$ws.Range('E38').Value = '  -0.09%  '

$ws.Range('E39').Value = '  +1.67%  '

$ws.Range('E40').Value = '  -0.43%  '

$ws.Range('E41').Value = '  +2.17%  '

$ws.Range('D42').Value = '0.0650'
$ws.Range('E42').Value = '  +7.60%  '

$ws.Range('D43').Value = '16.27'
$ws.Range('E43').Value = '  +6.44%  '

$ws.Range('D44').Value = '89.98'
$ws.Range('E44').Value = '  -1.74%  '

$ws.Range('D45').Value = '1.337.33'
$ws.Range('E45').Value = '  -1.37%  '

$ws.Range('D46').Value = '49.29'
$ws.Range('E46').Value = '  +40.49%  '

$ws.Range('E47').Value = '  +1.52%  '

$ws.Range('D48').Value = '2.80'
$ws.Range('E48').Value = '  +3.00%  '

$ws.Range('E49').Value = '  -1.04%  '

$ws.Range('D50').Value = '6.51'
$ws.Range('E50').Value = '  -2.26%  '

$ws.Range('D51').Value = '2.096.11'
$ws.Range('E51').Value = '  +1.56%  '

# Restore the default "Normal" style on the cells we forced to Text format,
# so only the value changes (no stray style/format diff is introduced).
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
